$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header row cells in-place:
#   columns A:J   "<Name>_old" -> "<Name>_FV2304"
#   column  K     "diff"        (unchanged)
#   columns L:U   "<Name>_new" -> "<Name>_FV2310"
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cur = [string]$cell.Value2
    $cell.Value = ($cur -replace "_old$", "_FV2304")
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cur = [string]$cell.Value2
    $cell.Value = ($cur -replace "_new$", "_FV2310")
}

# Add a table over A1:U64
$range = $ws.Range("A1:U64")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"

# Freeze panes on row 1
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
